{"js": "// Replace the whole body content (a Word/SDET password-policy test-case\n// write-up) with a Python unittest-style code block, per the commit:\n// \"Update ... from SDET sync (code block from Test Lead)\".\n//\n// Target paragraph texts, in order. `null` marks a blank paragraph.\nconst targetTexts = [\"import unittest\", null, null, null, \"    def test_less_than_8_chars(self):\", \"        password = \\\"abc1$\\\"\", null, \"    \", \"    def test_no_number(self):\", \"        password = \\\"abcd@xyz\\\" \", null, null, \"    def test_no_special_char(self):\", \"        password = \\\"abcd1234\\\"\", null, null, \"    def test_valid_8_chars(self):\", \"        password = \\\"abc1@def\\\"\", null, null, \"    def test_valid_gt_8_chars(self):\", \"        password = \\\"MyPass123!\\\"\", null, null, \"    def test_valid_new1(self):\", \"        password = \\\"1234@5678\\\"\", null, null, \"    def test_invalid_new2(self):\", \"        password = \\\"abcdefgh\\\"\", null, null, \"    def test_invalid_new3(self):\", \"        password = \\\"abcd1234\\\" \", null, null, \"    def test_valid_new4(self):\", \"        password = \\\"Ab1$xyz9\\\"\", null, null, \"def is_valid_password(password):\", \"    # Password validation logic\", \"    return True\", null, \"if __name__ == '__main__':\"];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\n\n// The very last paragraph in the body owns the section mark and can never\n// be fully removed, so keep paragraph[0] and paragraph[count-1] as anchors\n// and delete everything strictly between them.\nfor (let i = count - 2; i >= 1; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\nconst anchors = body.paragraphs;\nanchors.load(\"items\");\nawait context.sync();\n\nconst firstPara = anchors.items[0];\nconst lastPara = anchors.items[1];\n\n// Clear first so any leftover `xml:space=\"preserve\"` from the old text\n// doesn't linger on text that no longer needs it.\nfirstPara.clear();\nlastPara.clear();\n\n// First target line goes into the retained first paragraph.\nfirstPara.insertText(targetTexts[0], \"Start\");\n\n// Middle target lines get inserted (in order) right before the retained\n// last paragraph.\nconst middleTexts = targetTexts.slice(1, targetTexts.length - 1);\nfor (const text of middleTexts) {\n  lastPara.insertParagraph(text === null ? \"\" : text, \"Before\");\n}\n\n// Final target line goes into the retained last paragraph.\nlastPara.insertText(targetTexts[targetTexts.length - 1], \"Start\");\n\nawait context.sync();\n", "ps1": "# Replace the whole body content (a Word/SDET password-policy test-case\n# write-up) with a Python unittest-style code block, per the commit:\n# \"Update ... from SDET sync (code block from Test Lead)\".\n#\n# Target paragraph texts, in order. Empty strings stand for blank paragraphs.\n$targetTexts = @(\n    \"import unittest\",\n    \"\",\n    \"\",\n    \"\",\n    \"    def test_less_than_8_chars(self):\",\n    \"        password = `\"abc1`$`\"\",\n    \"\",\n    \"    \",\n    \"    def test_no_number(self):\",\n    \"        password = `\"abcd@xyz`\" \",\n    \"\",\n    \"\",\n    \"    def test_no_special_char(self):\",\n    \"        password = `\"abcd1234`\"\",\n    \"\",\n    \"\",\n    \"    def test_valid_8_chars(self):\",\n    \"        password = `\"abc1@def`\"\",\n    \"\",\n    \"\",\n    \"    def test_valid_gt_8_chars(self):\",\n    \"        password = `\"MyPass123!`\"\",\n    \"\",\n    \"\",\n    \"    def test_valid_new1(self):\",\n    \"        password = `\"1234@5678`\"\",\n    \"\",\n    \"\",\n    \"    def test_invalid_new2(self):\",\n    \"        password = `\"abcdefgh`\"\",\n    \"\",\n    \"\",\n    \"    def test_invalid_new3(self):\",\n    \"        password = `\"abcd1234`\" \",\n    \"\",\n    \"\",\n    \"    def test_valid_new4(self):\",\n    \"        password = `\"Ab1`$xyz9`\"\",\n    \"\",\n    \"\",\n    \"def is_valid_password(password):\",\n    \"    # Password validation logic\",\n    \"    return True\",\n    \"\",\n    \"if __name__ == '__main__':\"\n)\n\n$d = $word.ActiveDocument\n\n# The final paragraph in the body owns the section mark and can never be\n# fully removed, so keep paragraph 1 and the last paragraph as anchors and\n# delete everything strictly between them first.\n$count = $d.Paragraphs.Count\nfor ($i = $count - 1; $i -ge 2; $i--) {\n    $d.Paragraphs($i).Range.Delete()\n}\n\n# Set the first paragraph to the first target line.\n$d.Paragraphs(1).Range.Text = $targetTexts[0]\n\n# Insert a fresh blank paragraph right before the old last paragraph, then\n# delete that old last paragraph (which shifts down by one once the blank\n# paragraph is inserted ahead of it). This avoids inheriting any leftover\n# xml:space=\"preserve\" stickiness from the original text and gives us a\n# brand-new, \"clean\" last paragraph to build on.\n$oldLastIdx = $d.Paragraphs.Count\n$d.Paragraphs($oldLastIdx).Range.InsertParagraphBefore()\n$d.Paragraphs($oldLastIdx + 1).Range.Delete()\n\n# Insert all of the middle target lines (in order), each one right before\n# the current last (blank) paragraph, by re-querying the last index fresh\n# every time so we never hold a stale reference.\n$middleTexts = $targetTexts[1..($targetTexts.Length - 2)]\nforeach ($t in $middleTexts) {\n    $lastIdx = $d.Paragraphs.Count\n    $d.Paragraphs($lastIdx).Range.InsertParagraphBefore()\n    $d.Paragraphs($lastIdx).Range.Text = $t\n}\n\n# Finally, set the trailing blank paragraph to the last target line.\n$finalIdx = $d.Paragraphs.Count\n$d.Paragraphs($finalIdx).Range.Text = $targetTexts[$targetTexts.Length - 1]\n"}
